$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "lost_progress"
$ws.Range("B16").Value = "A key element of Foddian games is the ability to fall down and lose progress.  Designing levels around fish flopping proved difficult to create that kind of verticality, but, as you’ve just demonstrated, it is still quite easy to lose progress on a horizontal surface."

$ws.Range("B17").Select()
